# Refresh the cryptos snapshot: updated Price / Volume(1h) figures, and
# three coins reshuffled position in the ranking (rows 20/21 and 48-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.759.13"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "3.753.69"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "618.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("D7").Value = "3.751.07"
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.04%  "
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "4.372.97"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "3.750.25"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "69.796.92"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +21.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +5.70%  "
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "421.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "3.009.58"
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("E51").Value = "  +0.24%  "
